$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.043.01"
$ws.Range("E2").Value = "  -3.46%  "
$ws.Range("D3").Value = "3.506.09"
$ws.Range("E3").Value = "  -5.48%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.54"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.29"
$ws.Range("E6").Value = "  -4.65%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.607"
$ws.Range("E7").Value = "  -1.29%  "
$ws.Range("D8").Value = "3.500.67"
$ws.Range("E8").Value = "  -5.40%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.188"
$ws.Range("E10").Value = "  -6.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.63"
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.584"
$ws.Range("E12").Value = "  -4.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.35"
$ws.Range("E13").Value = "  -6.20%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000272"
$ws.Range("E14").Value = "  -5.43%  "
$ws.Range("D15").Value = "4.066.64"
$ws.Range("E15").Value = "  -5.63%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.52"
$ws.Range("E16").Value = "  -5.96%  "
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "627.91"
$ws.Range("E17").Value = "  -7.53%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.513.38"
$ws.Range("E18").Value = "  -5.53%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "69.095.84"
$ws.Range("E19").Value = "  -3.52%  "
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.30"
$ws.Range("E21").Value = "  -4.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.09"
$ws.Range("E22").Value = "  -4.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.887"
$ws.Range("E23").Value = "  -6.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.92"
$ws.Range("E24").Value = "  -9.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.25"
$ws.Range("E25").Value = "  -5.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.81"
$ws.Range("E26").Value = "  -4.77%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.63"
$ws.Range("E28").Value = "  -7.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.31"
$ws.Range("E29").Value = "  -11.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.62"
$ws.Range("E30").Value = "  -8.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.17"
$ws.Range("E31").Value = "  -8.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.52"
$ws.Range("E32").Value = "  -7.48%  "
$ws.Range("E33").Value = "  -9.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.02"
$ws.Range("E34").Value = "  -4.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "621.97"
$ws.Range("E35").Value = "  +5.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.76"
$ws.Range("E36").Value = "  -4.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.51"
$ws.Range("E37").Value = "  -14.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.103"
$ws.Range("E38").Value = "  -5.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "56.75"
$ws.Range("E39").Value = "  -4.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("E41").Value = "  -2.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.137"
$ws.Range("E42").Value = "  -6.63%  "
$ws.Range("D43").Value = "3.358.00"
$ws.Range("E43").Value = "  -8.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.328"
$ws.Range("E44").Value = "  -6.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.80"
$ws.Range("E45").Value = "  -8.14%  "
$ws.Range("D46").Value = "0.0₃0695"
$ws.Range("E46").Value = "  -10.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.56"
$ws.Range("E47").Value = "  -8.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.78"
$ws.Range("E48").Value = "  -3.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.130"
$ws.Range("E49").Value = "  -2.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.02"
$ws.Range("E50").Value = "  -2.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.65"
$ws.Range("E51").Value = "  +14.70%  "
